$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 503, shifting existing rows 503:632 down to 504:633
$ws.Rows("503:503").Insert()

# Populate the newly inserted row 503 with the new data record
$ws.Range("A503").Value = 9
$ws.Range("B503").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C503").Value = "Metropolitana"
$ws.Range("D503").Value = 45204
$ws.Range("D503").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E503").Value = 13
$ws.Range("F503").Value = 100112039
$ws.Range("G503").Value = "Ciboulette"
$ws.Range("H503").Value = "Sin especificar"
$ws.Range("I503").Value = "Primera"
$ws.Range("J503").Value = 430
$ws.Range("K503").Value = 1000
$ws.Range("L503").Value = 1200
$ws.Range("M503").Value = 1100
$ws.Range("N503").Value = "$/docena de atados"
$ws.Range("O503").Value = "Región Metropolitana"
$ws.Range("P503").Value = 367
$ws.Range("Q503").Value = 3
$ws.Range("R503").Value = "Hortaliza"
